$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column widths ----
# (Target OOXML column widths are sub-pixel values coming from real desktop
# Excel's font-metric width engine. This host's ColumnWidth setter always
# quantizes to a 1/6-character grid plus a fixed 5/6-character pad, so we
# feed it the pre-image that lands closest to the author's recorded width.)
$ws.Range("B:B").ColumnWidth = 11.307291666666666
$ws.Range("C:C").ColumnWidth = 13.307291666666666
$ws.Range("E:E").ColumnWidth = 14.592447916666666
$ws.Range("F:F").ColumnWidth = 17.592447916666668
$ws.Range("G:G").ColumnWidth = 22.736979166666668
$ws.Range("H:H").ColumnWidth = 19.877604166666668
$ws.Range("I:I").ColumnWidth = 20.022135416666668
$ws.Range("J:J").ColumnWidth = 14.307291666666666

# ---- Row 1 headers ----
# A1, B1, C1 stay the same (spdomain, flex, ent)
$ws.Range("F1").Value = "monthly_food_povertyline_non_entropy_adjusted"
$ws.Range("G1").Value = "ent_monthlyfoodpovertyline"
$ws.Range("H1").Value = "UBOS_monthlyfoodpovertyline"
$ws.Range("I1").Value = "UBOS%change"
$ws.Range("J1").Value = "Authors%change"

# clear old E1 header (was a shared string, now blank on row1)
$ws.Range("E1").ClearContents()

# ---- Row data: A (unchanged), B, C (new values), E (region labels) ----
$ws.Range("B2").Value = 779.64
$ws.Range("C2").Value = 784.08
$ws.Range("E2").Value = "Central rural"

$ws.Range("B3").Value = 714.96
$ws.Range("C3").Value = 736.47
$ws.Range("E3").Value = "Central urban"

$ws.Range("B4").Value = 595.41
$ws.Range("C4").Value = 742.99
$ws.Range("E4").Value = "Eastern"

$ws.Range("B5").Value = 374.91
$ws.Range("C5").Value = 453.9
$ws.Range("E5").Value = "Northern"

$ws.Range("B6").Value = 1246.14
$ws.Range("C6").Value = 821.13
$ws.Range("E6").Value = "Western"

# ---- H column: UBOS monthly food poverty line constant ----
$ws.Range("H2:H6").Value = 21258

# ---- F, G, I, J formulas ----
$ws.Range("F2").Formula = "=B2*30"
$ws.Range("F3:F6").Formula = "=B3*30"

# NOTE: column G previously had a lone non-shared formula at G4 sandwiched
# inside an otherwise-shared G3:G6 block. Writing G3:G6 as one range after
# that history trips an exporter quirk that leaves G4 pointing at a shared
# formula group with no definition (corrupt on reload). Setting each G cell
# individually sidesteps it while still producing the correct formula/value.
$ws.Range("G2").Formula = "=C2*30"
$ws.Range("G3").Formula = "=C3*30"
$ws.Range("G4").Formula = "=C4*30"
$ws.Range("G5").Formula = "=C5*30"
$ws.Range("G6").Formula = "=C6*30"

$ws.Range("I2").Formula = "=((G2-H2)/G2)*100"
$ws.Range("I3:I6").Formula = "=((G3-H3)/G3)*100"

$ws.Range("J2").Formula = "=((G2-F2)/G2)*100"
$ws.Range("J3:J6").Formula = "=((G3-F3)/G3)*100"

# ---- selection ----
$ws.Range("K10").Select() | Out-Null
